$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Remove the comment on G1 (also drops comments1.xml / vmlDrawing ref / legacyDrawing).
$ws.Range("G1").Comment.Delete()

# 2) Clear the old sample rows (3-6) completely: plain (unstyled) cells disappear on save,
#    while D (hyperlink style) and F (date style) keep their formatting shell but go empty.
$ws.Range("A3:A6").ClearContents()
$ws.Range("B3:C6").ClearContents()
$ws.Range("D3:D6").ClearContents()
$ws.Range("F3:F6").ClearContents()
$ws.Range("G3:G6").ClearContents()

# 3) Remove every existing hyperlink; we'll re-add the single one we keep (D2) further down.
$ws.Hyperlinks.Delete()

# 4) Re-populate row 2 with the new single user's data. Order matters: it controls the
#    order new shared strings are appended in sharedStrings.xml.
$ws.Range("D2").Value = "petya-pervyy-1999@mail.ru"
$ws.Range("B2").Value = "Петр"
$ws.Range("C2").Value = "Жигулёвский"
$ws.Range("E2").Value = "nWE#w(Qb"

# 5) New column F header + row2 value ("password" / "email password" pair).
$ws.Range("F1").Value = "Пароль от электронной почты"
$ws.Range("F2").Value = "ntvyjnf123"

# 6) Slide the old F column (birth date) / G column (gender) over to G / H for row 1 and 2.
$ws.Range("G1").Value = "Дата рождения"
$ws.Range("H1").Value = "Пол"

$ws.Range("G2").Value = 32874
$ws.Range("G2").NumberFormat = "m/d/yy"
$ws.Range("H2").Value = 2

# 7) Hyperlink for the email address cell.
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:petya-pervyy-1999@mail.ru") | Out-Null

# 8) Column widths: F grows to fit the new header, G is the (new) date column.
$ws.Columns.Item(6).ColumnWidth = 32.166666666666664
$ws.Columns.Item(7).ColumnWidth = 15.333333333333334

# 9) Selection as left by the edit.
$ws.Range("E3:E6").Select()
$excel.ActiveWindow.RangeSelection.Item(1).Activate() | Out-Null
